$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - 想去人数 (F) column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7480
$ws1.Range("F4").Value = 281
$ws1.Range("F5").Value = 5
$ws1.Range("F6").Value = 446
$ws1.Range("F7").Value = 4060
$ws1.Range("F8").Value = 322
$ws1.Range("F9").Value = 571
$ws1.Range("F10").Value = 274
$ws1.Range("F11").Value = 650
$ws1.Range("F12").Value = 130

# Sheet "全部类型" (sheet4) - 想去人数 (F) column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7481
$ws4.Range("F6").Value = 281
$ws4.Range("F7").Value = 5
$ws4.Range("F8").Value = 446
$ws4.Range("F9").Value = 4060
$ws4.Range("F10").Value = 322
$ws4.Range("F11").Value = 571
$ws4.Range("F12").Value = 274
$ws4.Range("F13").Value = 651
$ws4.Range("F15").Value = 130
